$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing B:C table one column to the left (B->A, C->B),
# preserving values/styles/column widths exactly (matches the diff,
# which moves the whole "oznacenie povodu / typ dokladu" table from
# columns B/C to columns A/B).
$ws.Columns("A").Delete()

# Row 8 already exists (previously the blank "C8" cell, now shifted to
# B8) - give it its real data.
$ws.Range("A8").Value = "Dobropis"
$ws.Range("B8").Value = "Faktúra - Dobropis"

# New row 9 with the second pair of added values.
$ws.Range("A9").Value = "Dávka platobných príkazov"
$ws.Range("B9").Value = "Vrátenie platby"

# Selection moves to A12 in the edited file.
$ws.Range("A12").Select()
